# Update the date line and the multiplication expressions in the table,
# per the commit: "Update master to output generated at c8c62b6"

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Header date
Replace-Text "2025-10-01 Wednesday" "2025-10-02 Thursday"

# Row 1
Replace-Text "316×3=" "377×4="
Replace-Text "498×8=" "754×8="
Replace-Text "119×8=" "830×2="
Replace-Text "268×7=" "429×3="
Replace-Text "751×5=" "905×7="

# Row 2
Replace-Text "501×6=" "151×9="
Replace-Text "192×8=" "223×9="
Replace-Text "418×8=" "222×4="
Replace-Text "370×5=" "951×5="
Replace-Text "796×7=" "188×7="

# Row 3
Replace-Text "746×9=" "726×2="
Replace-Text "348×5=" "381×5="
Replace-Text "975×4=" "228×9="
Replace-Text "638×9=" "735×9="
Replace-Text "694×2=" "309×3="

# Row 4
Replace-Text "462×3=" "304×2="
Replace-Text "923×5=" "423×8="
Replace-Text "186×2=" "838×7="
Replace-Text "988×2=" "832×5="
Replace-Text "482×3=" "506×4="

# Row 5
Replace-Text "782×3=" "434×7="
Replace-Text "349×9=" "488×2="
Replace-Text "502×8=" "356×7="
Replace-Text "439×9=" "552×5="
Replace-Text "894×6=" "894×4="
